$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7661
$ws.Range("K3").Value = 7921
$ws.Range("J4").Value = 1825
$ws.Range("K5").Value = 568
$ws.Range("K6").Value = 8822
$ws.Range("J7").Value = 28773
$ws.Range("K7").Value = 26636

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 486
$ws.Range("K3").Value = 526
$ws.Range("K7").Value = 1739

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 395
$ws.Range("K7").Value = 1119

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 152
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 446

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 203
$ws.Range("K7").Value = 623

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 190
$ws.Range("K6").Value = 108
$ws.Range("K7").Value = 448

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K4").Value = 97
$ws.Range("K7").Value = 797
$ws.Range("K8").Value = 1739
$ws.Range("K11").Value = 469
$ws.Range("K15").Value = 275
$ws.Range("K20").Value = 653
$ws.Range("K27").Value = 256
$ws.Range("K29").Value = 1468
$ws.Range("K30").Value = 106
$ws.Range("K31").Value = 315
$ws.Range("K33").Value = 1119
$ws.Range("K36").Value = 343
$ws.Range("K41").Value = 178
$ws.Range("K42").Value = 984
$ws.Range("K43").Value = 223
$ws.Range("K47").Value = 181
$ws.Range("K52").Value = 686
$ws.Range("K54").Value = 525
$ws.Range("K57").Value = 110
$ws.Range("K64").Value = 159
$ws.Range("K65").Value = 623
$ws.Range("K70").Value = 47
$ws.Range("K75").Value = 87
$ws.Range("K76").Value = 369
$ws.Range("J77").Value = 205
$ws.Range("K77").Value = 173
$ws.Range("K78").Value = 333
$ws.Range("K84").Value = 217
$ws.Range("K85").Value = 1228
$ws.Range("K89").Value = 397
$ws.Range("K92").Value = 98
$ws.Range("K94").Value = 360
$ws.Range("K95").Value = 446
$ws.Range("K99").Value = 448
$ws.Range("J101").Value = 28773
$ws.Range("K101").Value = 26636

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 123
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 81
$ws.Range("K3").Value = 120
$ws.Range("K7").Value = 525

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 410
$ws.Range("K6").Value = 435
$ws.Range("K7").Value = 1468

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K4").Value = 31
$ws.Range("K6").Value = 183
$ws.Range("K7").Value = 369

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 262
$ws.Range("K3").Value = 288
$ws.Range("K7").Value = 984

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 85
$ws.Range("K4").Value = 33
$ws.Range("K6").Value = 108
$ws.Range("K7").Value = 333

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 204
$ws.Range("K7").Value = 653

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 109
$ws.Range("K7").Value = 343

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 262
$ws.Range("K4").Value = 30
$ws.Range("K6").Value = 223
$ws.Range("K7").Value = 797

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 85
$ws.Range("K7").Value = 360

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 161
$ws.Range("K7").Value = 469

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K4").Value = 45
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 397

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 223

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 405
$ws.Range("K3").Value = 426
$ws.Range("K5").Value = 35
$ws.Range("K6").Value = 301
$ws.Range("K7").Value = 1228

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 72
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 205
$ws.Range("K7").Value = 173

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 190
$ws.Range("K7").Value = 686

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 206
$ws.Range("K4").Value = 42

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 159
